$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.374.45"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.158.70"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.74"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.03"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.395"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0860"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.00"
$ws.Range("E12").Value = "  +6.92%  "
$ws.Range("D13").Value = "2.479.48"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.11"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.814"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "2.155.90"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "39.338.59"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.20"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.30"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.75"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.69"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +8.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.65"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.08"
$ws.Range("E35").Value = "  +9.53%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.15"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.03"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "1.537.55"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +6.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0930"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +7.02%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.76"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.20"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "2.363.04"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("E51").Value = "  +0.33%  "
